# Apply cryptos list update (values scraped on Sat May 13 09:41:29 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.999.72'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.37%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.819.32'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.04%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.45%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.21'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.49%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.007'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.53%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4306'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.19%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3692'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.67%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07276'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.83%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8701'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.04%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '2.113.15'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +17.49%  '

$ws.Range("E12").Value = '  +5.24%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.420'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.36%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.636'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.29%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06987'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.75%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '81.21'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.56%  '

$ws.Range("E17").Value = '  +0.93%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008915'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.68%  '

$ws.Range("E19").Value = '  +0.66%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.27'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.57%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.062.27'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.56%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.209'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.42%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.04'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.70%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.350.06'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +16.37%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.34'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.06%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.887'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.29%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.43'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.62%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.242'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.925'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +13.13%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.03'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.84%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08983'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.65%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.186'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.67%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7518'
$ws.Range("D33").Style = "Normal"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.431'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.25%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.814'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.87%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.008'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.60%  '

$ws.Range("E37").Value = '  +4.97%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05250'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.42%  '

$ws.Range("E39").Value = '  +1.91%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5131'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.58%  '

$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1655'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.82%  '

$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.749'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +8.88%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.504'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.60%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.341'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.65%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '107.28'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.24%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.40'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.19%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.007'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.60%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.654'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.87%  '

$ws.Range("E49").Value = '  +2.14%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06225'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.66%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.849'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.96%  '
